$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B14: title change
$ws.Range("B14").Value = "Demo M4M Vocabulary"

# B17: clear creator ORCID value
$ws.Range("B17").Value = ""

# B21: modified datetime change
$ws.Range("B21").Value = "2023-06-08T05:50:35+00:00"

# Row 23: swap C23 (skos:altLabel...) and F23 (skos:broader...)
$ws.Range("C23").Value = "skos:broader(lookupColumn=""skos:prefLabel"" separator="","")"
$ws.Range("F23").Value = "skos:altLabel(separator="","")"

# Row 25: move "M4M26 vocabulary" from F25 to C25
$ws.Range("C25").Value = "M4M26 vocabulary"
$ws.Range("F25").Value = ""

# Row 26: move "M4M26 vocabulary" from F26 to C26
$ws.Range("C26").Value = "M4M26 vocabulary"
$ws.Range("F26").Value = ""

# Delete rows 27-35 entirely (shift up)
$ws.Range("A27:AM35").EntireRow.Delete()
